# Update schedule data (Línea 141) across all three sheets to the latest
# scrape snapshot (Última actualización 06:52:41 -> 07:23:38), growing each
# table with newly scraped rows and refreshed "Minutos" countdowns.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,1).Value = "Última actualización: 07:23:38"
$ws1.Cells.Item(3,1).Value = "Total filas: 55"
$ws1.Cells.Item(33,1).Value = "07:23:38"
$ws1.Cells.Item(33,2).Value = "07:23"
$ws1.Cells.Item(33,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(33,4).Value = 0
$ws1.Cells.Item(33,5).Value = "LP1912"
$ws1.Cells.Item(34,1).Value = "06:52:41"
$ws1.Cells.Item(34,2).Value = "07:23"
$ws1.Cells.Item(34,3).Value = "10_OLMOS"
$ws1.Cells.Item(34,4).Value = 31
$ws1.Cells.Item(34,5).Value = "LP1912"
$ws1.Cells.Item(35,1).Value = "07:23:38"
$ws1.Cells.Item(35,2).Value = "07:31"
$ws1.Cells.Item(35,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(35,4).Value = 8
$ws1.Cells.Item(35,5).Value = "LP1912"
$ws1.Cells.Item(36,1).Value = "07:23:38"
$ws1.Cells.Item(36,2).Value = "07:31"
$ws1.Cells.Item(36,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(36,4).Value = 8
$ws1.Cells.Item(36,5).Value = "LP1912"
$ws1.Cells.Item(37,1).Value = "07:23:38"
$ws1.Cells.Item(37,2).Value = "07:32"
$ws1.Cells.Item(37,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(37,4).Value = 9
$ws1.Cells.Item(37,5).Value = "LP1912"
$ws1.Cells.Item(38,1).Value = "07:23:38"
$ws1.Cells.Item(38,2).Value = "07:36"
$ws1.Cells.Item(38,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(38,4).Value = 13
$ws1.Cells.Item(38,5).Value = "LP1912"
$ws1.Cells.Item(39,1).Value = "07:23:38"
$ws1.Cells.Item(39,2).Value = "07:36"
$ws1.Cells.Item(39,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(39,4).Value = 13
$ws1.Cells.Item(39,5).Value = "LP1912"
$ws1.Cells.Item(40,1).Value = "07:23:38"
$ws1.Cells.Item(40,2).Value = "07:39"
$ws1.Cells.Item(40,3).Value = "10_OLMOS"
$ws1.Cells.Item(40,4).Value = 16
$ws1.Cells.Item(40,5).Value = "LP1912"
$ws1.Cells.Item(41,1).Value = "07:23:38"
$ws1.Cells.Item(41,2).Value = "07:47"
$ws1.Cells.Item(41,3).Value = "14_ABASTO"
$ws1.Cells.Item(41,4).Value = 24
$ws1.Cells.Item(41,5).Value = "LP1912"
$ws1.Cells.Item(42,1).Value = "07:23:38"
$ws1.Cells.Item(42,2).Value = "07:51"
$ws1.Cells.Item(42,3).Value = "215D_EL PATO"
$ws1.Cells.Item(42,4).Value = 28
$ws1.Cells.Item(42,5).Value = "LP1912"
$ws1.Cells.Item(43,1).Value = "07:23:38"
$ws1.Cells.Item(43,2).Value = "07:55"
$ws1.Cells.Item(43,3).Value = "10_OLMOS"
$ws1.Cells.Item(43,4).Value = 32
$ws1.Cells.Item(43,5).Value = "LP1912"
$ws1.Cells.Item(44,1).Value = "07:23:38"
$ws1.Cells.Item(44,2).Value = "07:59"
$ws1.Cells.Item(44,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(44,4).Value = 36
$ws1.Cells.Item(44,5).Value = "LP1912"
$ws1.Cells.Item(45,1).Value = "07:23:38"
$ws1.Cells.Item(45,2).Value = "08:03"
$ws1.Cells.Item(45,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(45,4).Value = 40
$ws1.Cells.Item(45,5).Value = "LP1912"
$ws1.Cells.Item(46,1).Value = "06:52:41"
$ws1.Cells.Item(46,2).Value = "08:06"
$ws1.Cells.Item(46,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(46,4).Value = 74
$ws1.Cells.Item(46,5).Value = "LP1912"
$ws1.Cells.Item(47,1).Value = "07:23:38"
$ws1.Cells.Item(47,2).Value = "08:12"
$ws1.Cells.Item(47,3).Value = "15_ABASTO"
$ws1.Cells.Item(47,4).Value = 49
$ws1.Cells.Item(47,5).Value = "LP1912"
$ws1.Cells.Item(48,1).Value = "07:23:38"
$ws1.Cells.Item(48,2).Value = "08:21"
$ws1.Cells.Item(48,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(48,4).Value = 58
$ws1.Cells.Item(48,5).Value = "LP1912"
$ws1.Cells.Item(49,1).Value = "07:23:38"
$ws1.Cells.Item(49,2).Value = "08:22"
$ws1.Cells.Item(49,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(49,4).Value = 59
$ws1.Cells.Item(49,5).Value = "LP1912"
$ws1.Cells.Item(50,1).Value = "07:23:38"
$ws1.Cells.Item(50,2).Value = "08:23"
$ws1.Cells.Item(50,3).Value = "215B_EL PATO"
$ws1.Cells.Item(50,4).Value = 60
$ws1.Cells.Item(50,5).Value = "LP1912"
$ws1.Cells.Item(51,1).Value = "07:23:38"
$ws1.Cells.Item(51,2).Value = "08:27"
$ws1.Cells.Item(51,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(51,4).Value = 64
$ws1.Cells.Item(51,5).Value = "LP1912"
$ws1.Cells.Item(52,1).Value = "06:52:41"
$ws1.Cells.Item(52,2).Value = "08:42"
$ws1.Cells.Item(52,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(52,4).Value = 110
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "07:23:38"
$ws1.Cells.Item(53,2).Value = "08:43"
$ws1.Cells.Item(53,3).Value = "14_ABASTO"
$ws1.Cells.Item(53,4).Value = 80
$ws1.Cells.Item(53,5).Value = "LP1912"
$ws1.Cells.Item(54,1).Value = "07:23:38"
$ws1.Cells.Item(54,2).Value = "08:44"
$ws1.Cells.Item(54,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(54,4).Value = 81
$ws1.Cells.Item(54,5).Value = "LP1912"
$ws1.Cells.Item(55,1).Value = "07:23:38"
$ws1.Cells.Item(55,2).Value = "08:54"
$ws1.Cells.Item(55,3).Value = "17_ROMERO"
$ws1.Cells.Item(55,4).Value = 91
$ws1.Cells.Item(55,5).Value = "LP1912"
$ws1.Cells.Item(56,1).Value = "07:23:38"
$ws1.Cells.Item(56,2).Value = "09:01"
$ws1.Cells.Item(56,3).Value = "215A_EL PATO"
$ws1.Cells.Item(56,4).Value = 98
$ws1.Cells.Item(56,5).Value = "LP1912"
$ws1.Cells.Item(57,1).Value = "07:23:38"
$ws1.Cells.Item(57,2).Value = "09:10"
$ws1.Cells.Item(57,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(57,4).Value = 107
$ws1.Cells.Item(57,5).Value = "LP1912"
$ws1.Cells.Item(58,1).Value = "07:23:38"
$ws1.Cells.Item(58,2).Value = "09:16"
$ws1.Cells.Item(58,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(58,4).Value = 113
$ws1.Cells.Item(58,5).Value = "LP1912"
$ws1.Cells.Item(59,1).Value = "07:23:38"
$ws1.Cells.Item(59,2).Value = "09:21"
$ws1.Cells.Item(59,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(59,4).Value = 118
$ws1.Cells.Item(59,5).Value = "LP1912"
$ws1.Cells.Item(60,1).Value = "07:23:38"
$ws1.Cells.Item(60,2).Value = "09:22"
$ws1.Cells.Item(60,3).Value = "17_ROMERO"
$ws1.Cells.Item(60,4).Value = 119
$ws1.Cells.Item(60,5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,1).Value = "Última actualización: 07:23:38"
$ws2.Cells.Item(3,1).Value = "Total filas: 8"
$ws2.Cells.Item(11,1).Value = "07:23:38"
$ws2.Cells.Item(11,2).Value = "07:51"
$ws2.Cells.Item(11,3).Value = "215D_EL PATO"
$ws2.Cells.Item(11,4).Value = 28
$ws2.Cells.Item(11,5).Value = "LP1912"
$ws2.Cells.Item(12,1).Value = "07:23:38"
$ws2.Cells.Item(12,2).Value = "08:23"
$ws2.Cells.Item(12,3).Value = "215B_EL PATO"
$ws2.Cells.Item(12,4).Value = 60
$ws2.Cells.Item(12,5).Value = "LP1912"
$ws2.Cells.Item(13,1).Value = "07:23:38"
$ws2.Cells.Item(13,2).Value = "09:01"
$ws2.Cells.Item(13,3).Value = "215A_EL PATO"
$ws2.Cells.Item(13,4).Value = 98
$ws2.Cells.Item(13,5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,1).Value = "Última actualización: 07:23:38"
$ws3.Cells.Item(3,1).Value = "Total filas: 12"
$ws3.Cells.Item(12,1).Value = "07:23:38"
$ws3.Cells.Item(12,2).Value = "07:35"
$ws3.Cells.Item(12,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(12,4).Value = 12
$ws3.Cells.Item(12,5).Value = "L6173"
$ws3.Cells.Item(14,1).Value = "07:23:38"
$ws3.Cells.Item(14,2).Value = "08:09"
$ws3.Cells.Item(14,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(14,4).Value = 46
$ws3.Cells.Item(14,5).Value = "L6203"
$ws3.Cells.Item(15,1).Value = "06:52:41"
$ws3.Cells.Item(15,2).Value = "08:31"
$ws3.Cells.Item(15,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(15,4).Value = 99
$ws3.Cells.Item(15,5).Value = "L6173"
$ws3.Cells.Item(16,1).Value = "07:23:38"
$ws3.Cells.Item(16,2).Value = "08:35"
$ws3.Cells.Item(16,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(16,4).Value = 72
$ws3.Cells.Item(16,5).Value = "L6173"
$ws3.Cells.Item(17,1).Value = "07:23:38"
$ws3.Cells.Item(17,2).Value = "09:08"
$ws3.Cells.Item(17,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(17,4).Value = 105
$ws3.Cells.Item(17,5).Value = "L6203"
